$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Update H3/I3 and H4/I4 values per diff
$ws.Range("H3").Value = -36
$ws.Range("H4").Value = -48

# Set I3/I4 as literal text "16-Sep-2025" (not an auto-converted date), while
# keeping their existing cell style (s="3") intact. Direct .Value assignment
# of a date-like string triggers Excel's smart date parsing + a new
# number-formatted style, so instead stage the text in a scratch cell that's
# explicitly formatted as Text, then paste-special (values only) into the
# target cells so their original style/format is preserved.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "16-Sep-2025"
$scratch.Copy()
$ws.Range("I3").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("I4").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false
